$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply updated crypto market data values (GitHub Actions scheduled refresh)

$ws.Range('D2').Value = '27.954.48'
$ws.Range('E2').Value = '  -4.25%  '
$ws.Range('D3').Value = '1.740.44'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('E5').Value = '  -3.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5787'
$ws.Range('E6').Value = '  -3.36%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2743'
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '23.14'
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06629'
$ws.Range('E10').Value = '  -4.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07560'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '1.743.86'
$ws.Range('E12').Value = '  -4.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.708'
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6029'
$ws.Range('E14').Value = '  -3.70%  '
$ws.Range('D15').Value = '1.977.22'
$ws.Range('E15').Value = '  -4.58%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '74.60'
$ws.Range('E16').Value = '  -3.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008750'
$ws.Range('E17').Value = '  -10.54%  '
$ws.Range('D18').Value = '27.944.19'
$ws.Range('E18').Value = '  -3.66%  '
$ws.Range('E19').Value = '  -3.84%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '205.56'
$ws.Range('E21').Value = '  -4.64%  '
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.637'
$ws.Range('E23').Value = '  -2.75%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.24'
$ws.Range('E25').Value = '  -3.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.038'
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1235'
$ws.Range('E27').Value = '  -4.03%  '
$ws.Range('E28').Value = '  -1.74%  '
$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06205'
$ws.Range('E29').Value = '  -3.19%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.386'
$ws.Range('E30').Value = '  -2.85%  '
$ws.Range('E31').Value = '  -3.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.746'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.741'
$ws.Range('E33').Value = '  -0.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.677'
$ws.Range('E34').Value = '  -2.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.039'
$ws.Range('E35').Value = '  -4.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6403'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.440'
$ws.Range('E37').Value = '  -3.93%  '
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('E39').Value = '  -4.49%  '
$ws.Range('D40').Value = '1.123.84'
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.158'
$ws.Range('E41').Value = '  -6.55%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8750'
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.04'
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('E45').Value = '  -4.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.35'
$ws.Range('E46').Value = '  -4.19%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000109'
$ws.Range('E47').Value = '  -3.28%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.577'
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.260'
$ws.Range('E49').Value = '  -1.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05377'
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4413'
$ws.Range('E51').Value = '  -2.63%  '
